$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unused trailing rows (17 and 18) so the sheet
# dimension shrinks from A1:F18 to A1:F16
$ws.Range("A17:F18").Delete()

# Row 2
$ws.Range("B2").Value = "NSE:ATAM"
$ws.Range("C2").Value = "NSE:BAJFINANCE"
$ws.Range("D2").Value = "NSE:NHPC"
$ws.Range("E2").Value = "NSE:GLENMARK"
$ws.Range("F2").Value = "NSE:BANDHANBNK"

# Row 3
$ws.Range("B3").Value = "NSE:CENTURYPLY"
$ws.Range("C3").Value = "NSE:BHARTIARTL"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "NSE:HINDUNILVR"
$ws.Range("F3").Value = "NSE:COALINDIA"

# Row 4
$ws.Range("B4").Value = "NSE:ETHOSLTD"
$ws.Range("C4").Value = "NSE:DBSTOCKBRO"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "NSE:LAURUSLABS"
$ws.Range("F4").Value = "NSE:HUDCO"

# Row 5
$ws.Range("B5").Value = "NSE:GODREJIND"
$ws.Range("C5").Value = "NSE:FMGOETZE"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "NSE:MARUTI"
$ws.Range("F5").Value = "NSE:JSWENERGY"

# Row 6
$ws.Range("B6").Value = "NSE:HPIL"
$ws.Range("C6").Value = "NSE:GRINDWELL"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "NSE:RAMCOCEM"
$ws.Range("F6").Value = "NSE:NATIONALUM"

# Row 7
$ws.Range("B7").Value = "NSE:IGARASHI"
$ws.Range("C7").Value = "NSE:ICEMAKE"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = "NSE:NHPC"

# Row 8
$ws.Range("B8").Value = "NSE:JSWENERGY"
$ws.Range("C8").Value = "NSE:KESORAMIND"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = "NSE:SAIL"

# Row 9
$ws.Range("B9").Value = "NSE:KDDL"
$ws.Range("C9").Value = "NSE:LICNETFSEN"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""

# Row 10
$ws.Range("B10").Value = "NSE:LAOPALA"
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""

# Row 11
$ws.Range("B11").Value = "NSE:LGBBROSLTD"
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""

# Row 12
$ws.Range("B12").Value = "NSE:MANORG"
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""

# Row 13
$ws.Range("B13").Value = "NSE:MOLDTECH"
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""

# Row 14
$ws.Range("B14").Value = "NSE:OMAXAUTO"
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""

# Row 15
$ws.Range("B15").Value = "NSE:RAIN"
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""

# Row 16
$ws.Range("B16").Value = "NSE:SAIL"
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = ""
